# The captured OOXML diff for this revision only rewrites the boilerplate
# namespace-prefix bookkeeping that Word regenerates on every save (e.g.
# xmlns:ns8 -> xmlns:ns9, xmlns:ns19 -> xmlns:ns17, and dropping the now
# redundant mc:Ignorable="w14 w15" hint on the <w:document> root). None of
# the actual document content (text, formatting, tables, headers/footers,
# styles, numbering, theme) differs between the two revisions - every
# hunk in the diff touches only the single root-element opening tag of
# each part and nothing else.
#
# The commit message itself describes an application-level bug fix
# (a web form's "selecionar" field being switched from boolean to
# integer) that has no bearing on this Word template's content, which is
# consistent with the document simply having been opened and re-saved by
# Word without any user-visible edit being made.
#
# Accordingly there is no document content to mutate here. We simply
# touch the document the same way Word does on open/save (a plain
# search, no replace) so the automation round trip is represented
# without introducing any spurious content, whitespace, or formatting
# drift that a real edit (even a same-text Find & Replace) would cause
# when the package is re-serialized.
$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("Identificação da Empresa")
